# Calendar event display html reads properly
# Reorders several rows of the schedule table so entries read in the
# correct (chronological / presentation) order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2 and 3: swap Unit / Classroom / Lecturer ---------------------
$ws.Range("C2").Value = "MITS5507"
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = "Sammy"

$ws.Range("C3").Value = "MITS5002"
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = "Mitch"

# --- Rows 5, 6, 7: rotate Time / Unit / Classroom / Lecturer / Mode -----
# old row5 -> new row6, old row6 -> new row7, old row7 -> new row5
$ws.Range("B5").Value = "8:00 AM to 9:00 AM"
$ws.Range("C5").Value = "MITS4001"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Jim"
$ws.Range("F5").Value = "Online"

$ws.Range("B6").Value = "8:00 AM to 10:00 AM"
$ws.Range("C6").Value = "MITS5501"
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = "Lewis"
$ws.Range("F6").Value = "F2F"

$ws.Range("B7").Value = "8:00 AM to 10:00 AM"
$ws.Range("C7").Value = "MITS4003"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = "Tom"
$ws.Range("F7").Value = "F2F"

# --- Rows 15 and 16: swap Time / Unit / Classroom / Lecturer ------------
$ws.Range("B15").Value = "2:00 PM to 3:00 PM"
$ws.Range("C15").Value = "MITS6500"
$ws.Range("D15").Value = 21
$ws.Range("E15").Value = "Keno"

$ws.Range("B16").Value = "2:00 PM to 4:00 PM"
$ws.Range("C16").Value = "MITS5003"
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = "Jay"

$wb.Save()
